$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows before row 4 so old row4 (header) becomes row6, etc.
$ws.Range("A4:A5").EntireRow.Insert()

$ws.Range("A4").Value = "A step with 0 as the number of times executed means it was skipped"

# Fix a few cells that changed content within the shifted block
$ws.Range("D11").Value = "3. TextInput(description)(some random text)(1)"
$ws.Range("B14").Value = "1. TextFile(description1)()(0)"
$ws.Range("C14").Value = "2. CsvFile(description2)()(0)"
$ws.Range("D14").Value = "3. Display()()(1)"
$ws.Range("B25").Value = "1. Output()()(1)"
$ws.Range("B26").Value = "1. Added one of each flow, 2 Numbers"
$ws.Range("C26").Value = "2. Entered some bad inputs where I could"
$ws.Range("D26").Value = "3. Ran flow analysis(details) option"
$ws.Range("E26").Value = "4. Deleted the flow"
$ws.Range("F26").Value = "5. Tried to run it again"

# New "goals" section starting at row 29
$ws.Range("A29").Value = "The goals of each step"
$ws.Range("A30").Value = 1
$ws.Range("B30").Value = "Check if Title, Text and TextInput creation / execution work correctly"
$ws.Range("A31").Value = 2
$ws.Range("B31").Value = "Check Number and Calculus Steps, if calculus can get the previous number steps as parameters"
$ws.Range("A32").Value = 3
$ws.Range("B32").Value = "See if TextFile and CsvFile work as expected and Display can receive them as parameters"
$ws.Range("A33").Value = 4
$ws.Range("B33").Value = "Check Calculus with no previous Number Steps"
$ws.Range("A34").Value = 5
$ws.Range("B34").Value = "Run Output for Title, Text and TextInput, to see if they are added to the file as expected"
$ws.Range("A35").Value = 6
$ws.Range("B35").Value = "Ran Display with no previous TextFile or CsvFile, and checked if Output adds Number and Calculus to the file as it should"
$ws.Range("A36").Value = 7
$ws.Range("B36").Value = "Checked if the flow analytis for Number work as they should"
$ws.Range("A37").Value = 8
$ws.Range("B37").Value = "Skipped both previous file Steps and then executed Display to see what would be available to be displayed "
$ws.Range("A38").Value = 9
$ws.Range("B38").Value = "Ran Calculus with only one previous Number step"
$ws.Range("A39").Value = 10
$ws.Range("B39").Value = "Checked if division by 0 si handled well"
$ws.Range("A40").Value = 11
$ws.Range("B40").Value = "Added 2 of each file type to see if Display can display any of them"
$ws.Range("A41").Value = 12
$ws.Range("B41").Value = "Wanted to see if a flow is deleted correctly"
$ws.Range("A42").Value = 13
$ws.Range("B42").Value = "Checked if TextFile and CsvFile are added to the Output file as expected"
$ws.Range("A43").Value = 14
$ws.Range("B43").Value = "Ran Calculus with no previous Number steps"
$ws.Range("A44").Value = 15
$ws.Range("B44").Value = "Ran Calculus with no previous Number steps"
$ws.Range("A45").Value = 16
$ws.Range("B45").Value = "Checked if Calculus works well with a lot of previous Numbers"
$ws.Range("A46").Value = 17
$ws.Range("B46").Value = "Executed the second Number step after the Calculus and then ran Output to check if I ge the expected result"
$ws.Range("A47").Value = 18
$ws.Range("B47").Value = "Entered bad input for more steps and ran flow analysis to check skips and errors"
$ws.Range("A48").Value = 19
$ws.Range("B48").Value = "Executed Output with no previous Steps"
$ws.Range("A49").Value = 20
$ws.Range("B49").Value = "Created a flow with each step and 2 Number Steps, executed it with some bad inputs, ran code analysis and then deleted the flow"

# Column width + selection updates (bestFit widths recomputed by Excel for the new content)
$ws.Columns("B").ColumnWidth = 33.1640625
$ws.Columns("C").ColumnWidth = 35.6640625

$ws.Range("B38,E33").Select()

